# Add two new "剑指 Offer" problems (58 - I and 58 - II) as new rows
# at the bottom of the tracking sheet, matching the existing table layout:
#   A = index number, B = problem title (hyperlinked), C = topic/tag,
#   D = difficulty, each new row styled like its neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 158: 剑指 Offer 58 - I. 翻转单词顺序 -------------------------------
$ws.Range("A158").Value = 159
$ws.Range("B158").Value = "剑指 Offer 58 - I. 翻转单词顺序"
$ws.Hyperlinks.Add($ws.Range("B158"), "https://leetcode.cn/problems/fan-zhuan-dan-ci-shun-xu-lcof/") | Out-Null
$ws.Range("B158").Style = "Hyperlink"
$ws.Range("B158").HorizontalAlignment = -4131
$ws.Range("C158").Value = "字符串"
$ws.Range("D158").Value = "简单"

# --- Row 159: 剑指 Offer 58 - II. 左旋转字符串 ------------------------------
$ws.Range("A159").Value = 160
$ws.Range("B159").Value = "剑指 Offer 58 - II. 左旋转字符串"
$ws.Hyperlinks.Add($ws.Range("B159"), "https://leetcode.cn/problems/zuo-xuan-zhuan-zi-fu-chuan-lcof/") | Out-Null
$ws.Range("B159").Style = "Hyperlink"
$ws.Range("B159").HorizontalAlignment = -4131
$ws.Range("C159").Value = "字符串"
$ws.Range("D159").Value = "简单"

# Move the visible selection / scroll position the way the author left it.
$win = $excel.ActiveWindow
$win.ScrollRow = 128
$win.ScrollColumn = 1
$ws.Range("E157").Select() | Out-Null
